$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update values: B3 403 -> 402, B7 398 -> 399
$ws.Range("B3").Value = 402
$ws.Range("B7").Value = 399

# Move the selection to B7 (matches the saved selection in the diff)
$ws.Range("B7").Select()
